$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5999.5
$ws.Range("I43").Value = 5999
$ws.Range("K43").Value = 5999
$ws.Range("M43").Value = -5930
$ws.Range("H70").Value = 1667.3334
$ws.Range("I70").Value = 1458.3334
$ws.Range("J70").Value = 1771.8334
$ws.Range("K70").Value = 4375.0002
$ws.Range("L70").Value = 5315.5002
$ws.Range("M70").Value = -4105.0002
$ws.Range("N70").Value = -5855.5002
$ws.Range("H73").Value = 1667.3334
$ws.Range("I73").Value = 1458.3334
$ws.Range("J73").Value = 1771.8334
$ws.Range("K73").Value = 4375.0002
$ws.Range("L73").Value = 5315.5002
$ws.Range("M73").Value = -3439.0002
$ws.Range("N73").Value = -7187.5002
$ws.Range("H88").Value = 4874.3335
$ws.Range("I88").Value = 5571.2856
$ws.Range("K88").Value = 5571.2856
$ws.Range("M88").Value = -5165.2856
$ws.Range("H91").Value = 4874.3335
$ws.Range("I91").Value = 5571.2856
$ws.Range("K91").Value = 5571.2856
$ws.Range("M91").Value = -4167.2856
$ws.Range("H106").Value = 5099.4287
$ws.Range("I106").Value = 5314.769
$ws.Range("K106").Value = 5314.769
$ws.Range("M106").Value = -4683.769
$ws.Range("H113").Value = 4289.5
$ws.Range("I113").Value = 4213.5713
$ws.Range("J113").Value = 4466.6665
$ws.Range("K113").Value = 4213.5713
$ws.Range("L113").Value = 4466.6665
$ws.Range("M113").Value = -959.5712999999996
$ws.Range("N113").Value = -10974.6665
$ws.Range("H116").Value = 50774.832
$ws.Range("I116").Value = 24875
$ws.Range("K116").Value = 24875
$ws.Range("M116").Value = -21433
$ws.Range("H118").Value = 1784.125
$ws.Range("I118").Value = 854.8
$ws.Range("K118").Value = 2564.4
$ws.Range("M118").Value = -907.3999999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 993.36365
$ws.Range("I2").Value = 992.8
$ws.Range("K2").Value = 992.8
$ws.Range("M2").Value = -879.8
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0
$ws.Range("H32").Value = 136940.81
$ws.Range("I32").Value = 145228.88
$ws.Range("K32").Value = 145228.88
$ws.Range("M32").Value = -144941.88
$ws.Range("H41").Value = 3414.5
$ws.Range("I41").Value = 1802.5
$ws.Range("J41").Value = 8250.5
$ws.Range("K41").Value = 1802.5
$ws.Range("L41").Value = 8250.5
$ws.Range("M41").Value = -1388.5
$ws.Range("N41").Value = -9078.5
$ws.Range("H45").Value = 2079.8333
$ws.Range("I45").Value = 1784.2222
$ws.Range("J45").Value = 2966.6667
$ws.Range("K45").Value = 1784.2222
$ws.Range("L45").Value = 2966.6667
$ws.Range("M45").Value = -1407.2222
$ws.Range("N45").Value = -3720.6667
$ws.Range("H116").Value = 993.36365
$ws.Range("I116").Value = 992.8
$ws.Range("K116").Value = 992.8
$ws.Range("M116").Value = 1301.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 993.36365
$ws.Range("I3").Value = 992.8
$ws.Range("K3").Value = 992.8
$ws.Range("M3").Value = -878.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 207400
$ws.Range("I16").Value = 9250
$ws.Range("K16").Value = 9250
$ws.Range("M16").Value = -8963
$ws.Range("H86").Value = 105514.78
$ws.Range("J86").Value = 9086.25
$ws.Range("L86").Value = 9086.25
$ws.Range("N86").Value = -11332.25
$ws.Range("H89").Value = 105514.78
$ws.Range("J89").Value = 9086.25
$ws.Range("L89").Value = 45431.25
$ws.Range("N89").Value = -56663.25
$ws.Range("H107").Value = 1576.12
$ws.Range("I107").Value = 1463.3846
$ws.Range("J107").Value = 1698.25
$ws.Range("K107").Value = 1463.3846
$ws.Range("L107").Value = 1698.25
$ws.Range("M107").Value = 456.6153999999999
$ws.Range("N107").Value = -5538.25
$ws.Range("H113").Value = 207400
$ws.Range("I113").Value = 9250
$ws.Range("K113").Value = 9250
$ws.Range("M113").Value = -7080
$ws.Range("H122").Value = 9147.583000000001
$ws.Range("I122").Value = 2123.7666
$ws.Range("J122").Value = 44266.668
$ws.Range("K122").Value = 6371.2998
$ws.Range("L122").Value = 132800.004
$ws.Range("M122").Value = -3921.2998
$ws.Range("N122").Value = -137700.004
$ws.Range("H132").Value = 3559.1052
$ws.Range("I132").Value = 3441.6
$ws.Range("J132").Value = 3999.75
$ws.Range("K132").Value = 10324.8
$ws.Range("L132").Value = 11999.25
$ws.Range("M132").Value = -7794.799999999999
$ws.Range("N132").Value = -17059.25
$ws.Range("H133").Value = 49143.6
$ws.Range("I133").Value = 45722
$ws.Range("J133").Value = 49999
$ws.Range("K133").Value = 45722
$ws.Range("L133").Value = 49999
$ws.Range("M133").Value = -43192
$ws.Range("N133").Value = -55059

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 114.181816
$ws.Range("I2").Value = 135.375
$ws.Range("K2").Value = 812.25
$ws.Range("M2").Value = -699.25
$ws.Range("H4").Value = 7202309
$ws.Range("I4").Value = 16364357
$ws.Range("K4").Value = 49093071
$ws.Range("M4").Value = -49092959
$ws.Range("H29").Value = 489
$ws.Range("J29").Value = 444
$ws.Range("L29").Value = 1332
$ws.Range("N29").Value = -1886
$ws.Range("H113").Value = 23789.346
$ws.Range("J113").Value = 32390.684
$ws.Range("L113").Value = 97172.052
$ws.Range("N113").Value = -101512.052
$ws.Range("H129").Value = 422363.66
$ws.Range("I129").Value = 1669366.1
$ws.Range("J129").Value = 6696.1665
$ws.Range("K129").Value = 5008098.300000001
$ws.Range("L129").Value = 20088.4995
$ws.Range("M129").Value = -5003098.300000001
$ws.Range("N129").Value = -30088.4995
$ws.Range("H131").Value = 6062.909
$ws.Range("I131").Value = 798
$ws.Range("K131").Value = 2394
$ws.Range("M131").Value = 2646
$ws.Range("H134").Value = 8653
$ws.Range("I134").Value = 5150
$ws.Range("J134").Value = 12656.429
$ws.Range("K134").Value = 15450
$ws.Range("L134").Value = 37969.287
$ws.Range("M134").Value = -10380
$ws.Range("N134").Value = -48109.287

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5008
$ws.Range("I102").Value = 5406.8
$ws.Range("K102").Value = 5406.8
$ws.Range("M102").Value = -3784.8
$ws.Range("H122").Value = 2385.3928
$ws.Range("I122").Value = 1710.8889
$ws.Range("J122").Value = 2704.8948
$ws.Range("K122").Value = 5132.6667
$ws.Range("L122").Value = 8114.6844
$ws.Range("M122").Value = -2682.6667
$ws.Range("N122").Value = -13014.6844
$ws.Range("H132").Value = 16321.24
$ws.Range("I132").Value = 16584.041
$ws.Range("J132").Value = 10014
$ws.Range("K132").Value = 49752.12300000001
$ws.Range("L132").Value = 30042
$ws.Range("M132").Value = -47222.12300000001
$ws.Range("N132").Value = -35102

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3831
$ws.Range("I122").Value = 3052.6667
$ws.Range("J122").Value = 4998.5
$ws.Range("K122").Value = 9158.000100000001
$ws.Range("L122").Value = 14995.5
$ws.Range("M122").Value = -6708.000100000001
$ws.Range("N122").Value = -19895.5
$ws.Range("H132").Value = 3208.1667
$ws.Range("I132").Value = 2543.1667
$ws.Range("J132").Value = 3540.6667
$ws.Range("K132").Value = 7629.500100000001
$ws.Range("L132").Value = 10622.0001
$ws.Range("M132").Value = -5099.500100000001
$ws.Range("N132").Value = -15682.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2414.818
$ws.Range("I113").Value = 1881.3846
$ws.Range("K113").Value = 5644.1538
$ws.Range("M113").Value = -3474.1538
